$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = "check chatgpt HIAN log for more info"
$ws.Range("B8").NumberFormat = "@"

$ws.Range("E11").Select()
